$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# --- Increase font size (Default -> Large) for engineTemp, controllerTemp, motorCurrent rows ---
$ws.Range("C7").Value = "Large"
$ws.Range("C8").Value = "Large"
$ws.Range("C9").Value = "Large"

# --- Put proper legends into dials: GB text "50" -> "100" for the two SingleUseId rows ---
# Leading apostrophe forces these digit-only strings to be stored as text (shared string),
# matching the original column's text-typed cells; reset the style afterwards so no new
# cell-level formatting (e.g. quote-prefix) is left behind.
$ws.Range("F10").Value = "'100"
$ws.Range("F11").Value = "'100"

# --- Add a new "tenths of km" partial-distance row (row 13 placeholder + row 14 partialDistance) ---
$ws.Range("B13").Value = "SingleUseId13"
$ws.Range("C13").Value = "Large"
$ws.Range("D13").Value = "Left"
$ws.Range("E13").Value = "LTR"
$ws.Range("F13").Value = "'0"

$ws.Range("B14").Value = "partialDistance"
$ws.Range("C14").Value = "Default"
$ws.Range("D14").Value = "Left"
$ws.Range("E14").Value = "LTR"
$ws.Range("F14").Value = "<number>"
$ws.Range("G14").Value = "Large"

# Reset styles on the cells where we used a quote-prefix to force text storage, so no stray
# cell-level number formatting is introduced.
$ws.Range("F10:F11").Style = "Normal"
$ws.Range("F13").Style = "Normal"
